# Apply the edit described by the diff:
#  - Insert a brand-new data row at sheet row 333 (pushing the former
#    rows 333..375 down by one).
#  - Insert a second brand-new data row at (the now-shifted) row 363
#    (pushing the remaining rows down by one more).
# Net effect: 2 new rows of data, old data otherwise unchanged but
# shifted down; dimension grows from A1:R375 to A1:R377.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Insert new row at 333 and populate it.
# ---------------------------------------------------------------------
$ws.Rows.Item(333).Insert()

$ws.Cells.Item(333, 1).Value = 5
$ws.Cells.Item(333, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(333, 3).Value = "Maule"
$ws.Cells.Item(333, 4).Value = 44748
$ws.Cells.Item(333, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(333, 5).Value = 7
$ws.Cells.Item(333, 6).Value = 100112032
$ws.Cells.Item(333, 7).Value = "Zapallo italiano"
$ws.Cells.Item(333, 8).Value = "Sin especificar"
$ws.Cells.Item(333, 9).Value = "Primera"
$ws.Cells.Item(333, 10).Value = 300
$ws.Cells.Item(333, 11).Value = 12000
$ws.Cells.Item(333, 12).Value = 12000
$ws.Cells.Item(333, 13).Value = 12000
$ws.Cells.Item(333, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(333, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(333, 16).Value = 240
$ws.Cells.Item(333, 17).Value = 50
$ws.Cells.Item(333, 18).Value = "Hortaliza"

# ---------------------------------------------------------------------
# 2) Insert a second new row at (post-shift) row 363 and populate it.
# ---------------------------------------------------------------------
$ws.Rows.Item(363).Insert()

$ws.Cells.Item(363, 1).Value = 5
$ws.Cells.Item(363, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(363, 3).Value = "Maule"
$ws.Cells.Item(363, 4).Value = 44747
$ws.Cells.Item(363, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(363, 5).Value = 7
$ws.Cells.Item(363, 6).Value = 100112032
$ws.Cells.Item(363, 7).Value = "Zapallo italiano"
$ws.Cells.Item(363, 8).Value = "Sin especificar"
$ws.Cells.Item(363, 9).Value = "Primera"
$ws.Cells.Item(363, 10).Value = 300
$ws.Cells.Item(363, 11).Value = 11000
$ws.Cells.Item(363, 12).Value = 11000
$ws.Cells.Item(363, 13).Value = 11000
$ws.Cells.Item(363, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(363, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(363, 16).Value = 220
$ws.Cells.Item(363, 17).Value = 50
$ws.Cells.Item(363, 18).Value = "Hortaliza"
